$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as text, preserving the cell's existing style
# (plain .Value assignment on a numeric-looking string like "231.89"
#  gets auto-coerced to a Double by Excel, which would change the
#  stored cell type from a string to a number - not what we want here,
#  since the source sheet keeps these "Price" figures as literal text).
function Set-TextValue($range, $value) {
    $origStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = $origStyle
}

# Refresh Price (D) / Volume(1h) (E) figures for this run's snapshot
Set-TextValue $ws.Range("D2") "43.924.03"
$ws.Range("E2").Value = "  +4.60%  "
Set-TextValue $ws.Range("D3") "2.278.75"
$ws.Range("E3").Value = "  +2.13%  "
$ws.Range("E4").Value = "  -0.06%  "
Set-TextValue $ws.Range("D5") "231.89"
$ws.Range("E5").Value = "  +0.21%  "
Set-TextValue $ws.Range("D6") "0.628"
$ws.Range("E6").Value = "  +0.66%  "
Set-TextValue $ws.Range("D7") "61.58"
$ws.Range("E7").Value = "  +1.11%  "
$ws.Range("E8").Value = "  -0.03%  "
Set-TextValue $ws.Range("D9") "0.423"
$ws.Range("E9").Value = "  +5.31%  "
Set-TextValue $ws.Range("D10") "0.0950"
$ws.Range("E10").Value = "  +6.84%  "
Set-TextValue $ws.Range("D11") "57.82"
$ws.Range("E11").Value = "  -1.92%  "
$ws.Range("E12").Value = "  +0.51%  "
Set-TextValue $ws.Range("D13") "2.616.98"
$ws.Range("E13").Value = "  +2.10%  "
Set-TextValue $ws.Range("D14") "15.81"
$ws.Range("E14").Value = "  +0.96%  "
Set-TextValue $ws.Range("D15") "23.84"
$ws.Range("E15").Value = "  +9.27%  "
Set-TextValue $ws.Range("D16") "5.84"
$ws.Range("E16").Value = "  +4.99%  "
Set-TextValue $ws.Range("D17") "0.813"
$ws.Range("E17").Value = "  +1.80%  "
Set-TextValue $ws.Range("D18") "2.287.06"
$ws.Range("E18").Value = "  +2.64%  "
Set-TextValue $ws.Range("D19") "43.757.46"
$ws.Range("E19").Value = "  +4.51%  "
Set-TextValue $ws.Range("D20") "0.0₃0942"
$ws.Range("E20").Value = "  +5.63%  "
Set-TextValue $ws.Range("D21") "73.29"
$ws.Range("E21").Value = "  +0.89%  "
Set-TextValue $ws.Range("D22") "6.25"
$ws.Range("E22").Value = "  +3.46%  "
Set-TextValue $ws.Range("D23") "251.26"
$ws.Range("E23").Value = "  +0.70%  "
$ws.Range("E24").Value = "  -0.06%  "
Set-TextValue $ws.Range("D25") "2.57"
$ws.Range("E25").Value = "  +7.65%  "
Set-TextValue $ws.Range("D26") "2.37"
$ws.Range("E26").Value = "  +2.40%  "
Set-TextValue $ws.Range("D27") "9.87"
$ws.Range("E27").Value = "  +2.75%  "
Set-TextValue $ws.Range("D28") "171.09"
$ws.Range("E28").Value = "  +2.11%  "
Set-TextValue $ws.Range("D29") "0.141"
$ws.Range("E29").Value = "  -0.88%  "
Set-TextValue $ws.Range("D30") "20.62"
$ws.Range("E30").Value = "  +3.50%  "
Set-TextValue $ws.Range("D31") "1.47"
$ws.Range("E31").Value = "  +4.87%  "
Set-TextValue $ws.Range("D32") "2.67"
$ws.Range("E32").Value = "  +1.60%  "
$ws.Range("E33").Value = "  +0.37%  "
Set-TextValue $ws.Range("D34") "4.80"
$ws.Range("E34").Value = "  +3.92%  "
Set-TextValue $ws.Range("D35") "5.05"
$ws.Range("E35").Value = "  +2.12%  "
Set-TextValue $ws.Range("D36") "0.0663"
$ws.Range("E36").Value = "  +5.57%  "
Set-TextValue $ws.Range("D37") "6.54"
$ws.Range("E37").Value = "  -1.75%  "
$ws.Range("E38").Value = "  +2.55%  "
Set-TextValue $ws.Range("D39") "3.63"
$ws.Range("E39").Value = "  -1.55%  "
$ws.Range("E40").Value = "  +4.45%  "
$ws.Range("E41").Value = "  +0.02%  "
$ws.Range("E42").Value = "  +2.50%  "
Set-TextValue $ws.Range("D43") "0.000224"
$ws.Range("E43").Value = "  -13.16%  "
Set-TextValue $ws.Range("D44") "4.52"
$ws.Range("E44").Value = "  -5.78%  "
Set-TextValue $ws.Range("D47") "98.28"
$ws.Range("E47").Value = "  -0.75%  "
Set-TextValue $ws.Range("D48") "1.472.93"
$ws.Range("E48").Value = "  +0.16%  "
Set-TextValue $ws.Range("D49") "16.79"
$ws.Range("E49").Value = "  +1.60%  "
$ws.Range("E50").Value = "  +1.34%  "
$ws.Range("E51").Value = "  -1.19%  "

# Rows 45/46 swapped: Cronos now ranks above TrustWalletToken
$ws.Range("B45").Value = "Cronos"
$ws.Range("C45").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue $ws.Range("D45") "0.0974"
$ws.Range("E45").Value = "  -0.58%  "
$ws.Range("B46").Value = "TrustWalletToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-TextValue $ws.Range("D46") "1.22"
$ws.Range("E46").Value = "  +0.73%  "
